$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=4;   D='[0, 0, 0, 0, 0, 0, 0]'; E="[]" },
    @{Row=6;   D='[1, 0, 0, 1, 0, 0, 0]'; E="['Normal', 'ParamViolation']" },
    @{Row=7;   D='[0, 0, 0, 0, 0, 0, 0]'; E="[]" },
    @{Row=11;  D='[1, 0, 0, 1, 0, 0, 0]'; E="['Normal', 'ParamViolation']" },
    @{Row=24;  D='[1, 0, 0, 0, 0, 0, 0]'; E="['Normal']" },
    @{Row=26;  D='[0, 0, 0, 0, 0, 0, 0]'; E="[]" },
    @{Row=29;  D='[1, 0, 0, 0, 0, 0, 1]'; E="['Normal', 'SoftwareFault']" },
    @{Row=35;  D='[1, 0, 1, 0, 0, 0, 0]'; E="['Normal', 'HardwareFault']" },
    @{Row=36;  D='[1, 1, 1, 0, 0, 0, 0]'; E="['Normal', 'SurroundingEnvironment', 'HardwareFault']" },
    @{Row=38;  D='[0, 0, 1, 0, 0, 0, 0]'; E="['HardwareFault']" },
    @{Row=39;  D='[1, 0, 1, 0, 0, 0, 1]'; E="['Normal', 'HardwareFault', 'SoftwareFault']" },
    @{Row=54;  D='[0, 0, 0, 0, 0, 1, 0]'; E="['CommunicationIssue']" },
    @{Row=56;  D='[0, 0, 0, 0, 0, 0, 0]'; E="[]" },
    @{Row=61;  D='[0, 0, 0, 0, 0, 0, 1]'; E="['SoftwareFault']" },
    @{Row=81;  D='[1, 0, 1, 0, 0, 0, 0]'; E="['Normal', 'HardwareFault']" },
    @{Row=84;  D='[1, 0, 0, 1, 0, 0, 0]'; E="['Normal', 'ParamViolation']" },
    @{Row=93;  D='[1, 0, 1, 0, 0, 0, 1]'; E="['Normal', 'HardwareFault', 'SoftwareFault']" },
    @{Row=109; D='[1, 1, 0, 0, 0, 0, 0]'; E="['Normal', 'SurroundingEnvironment']" },
    @{Row=113; D='[1, 0, 1, 0, 0, 0, 0]'; E="['Normal', 'HardwareFault']" },
    @{Row=116; D='[1, 0, 0, 0, 0, 0, 1]'; E="['Normal', 'SoftwareFault']" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
